$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 334.5
$ws.Range("I33").Value = 334.5
$ws.Range("K33").Value = 334.5
$ws.Range("M33").Value = -105.5
$ws.Range("H64").Value = 3999.4
$ws.Range("I64").Value = 3998.5
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3998.5
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3750.5
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3999.4
$ws.Range("I67").Value = 3998.5
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3998.5
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -3140.5
$ws.Range("N67").Value = -5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3075.84
$ws.Range("I32").Value = 3162.375
$ws.Range("K32").Value = 3162.375
$ws.Range("M32").Value = -2875.375
$ws.Range("H53").Value = 14999
$ws.Range("I53").Value = 14999
$ws.Range("K53").Value = 14999
$ws.Range("M53").Value = -14317
$ws.Range("H74").Value = 1761.1364
$ws.Range("I74").Value = 1457.4117
$ws.Range("J74").Value = 2793.8
$ws.Range("K74").Value = 1457.4117
$ws.Range("L74").Value = 2793.8
$ws.Range("M74").Value = -583.4117000000001
$ws.Range("N74").Value = -4541.8
$ws.Range("H77").Value = 1761.1364
$ws.Range("I77").Value = 1457.4117
$ws.Range("J77").Value = 2793.8
$ws.Range("K77").Value = 7287.058500000001
$ws.Range("L77").Value = 13969
$ws.Range("M77").Value = -2919.058500000001
$ws.Range("N77").Value = -22705
$ws.Range("H110").Value = 2462.1667
$ws.Range("I110").Value = 1854.6
$ws.Range("K110").Value = 1854.6
$ws.Range("M110").Value = 190.4000000000001
$ws.Range("H122").Value = 6195.607
$ws.Range("I122").Value = 6595.857
$ws.Range("J122").Value = 4994.857
$ws.Range("K122").Value = 19787.571
$ws.Range("L122").Value = 14984.571
$ws.Range("M122").Value = -17337.571
$ws.Range("N122").Value = -19884.571
$ws.Range("H132").Value = 1828.8125
$ws.Range("I132").Value = 1832.9286
$ws.Range("K132").Value = 5498.7858
$ws.Range("M132").Value = -2968.7858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2490.5557
$ws.Range("I99").Value = 2677.25
$ws.Range("K99").Value = 2677.25
$ws.Range("M99").Value = -1179.25
$ws.Range("H107").Value = 9518.5
$ws.Range("I107").Value = 9022.200000000001
$ws.Range("J107").Value = 12000
$ws.Range("K107").Value = 9022.200000000001
$ws.Range("L107").Value = 12000
$ws.Range("M107").Value = -7102.200000000001
$ws.Range("N107").Value = -15840
$ws.Range("H134").Value = 5964.636
$ws.Range("I134").Value = 3085.5715
$ws.Range("K134").Value = 9256.7145
$ws.Range("M134").Value = -6721.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8647.9
$ws.Range("I31").Value = 8310.75
$ws.Range("K31").Value = 8310.75
$ws.Range("M31").Value = -8015.75
$ws.Range("H34").Value = 8647.9
$ws.Range("I34").Value = 8310.75
$ws.Range("K34").Value = 8310.75
$ws.Range("M34").Value = -8108.75
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H58").Value = 2663.8333
$ws.Range("I58").Value = 2897.889
$ws.Range("J58").Value = 1961.6666
$ws.Range("K58").Value = 2897.889
$ws.Range("L58").Value = 1961.6666
$ws.Range("M58").Value = -2694.889
$ws.Range("N58").Value = -2367.6666
$ws.Range("H122").Value = 1000.1
$ws.Range("I122").Value = 1017.4286
$ws.Range("J122").Value = 959.6667
$ws.Range("K122").Value = 3052.2858
$ws.Range("L122").Value = 2879.0001
$ws.Range("M122").Value = -602.2857999999997
$ws.Range("N122").Value = -7779.0001
$ws.Range("H134").Value = 2472.2307
$ws.Range("I134").Value = 2416.889
$ws.Range("J134").Value = 2596.75
$ws.Range("K134").Value = 7250.667
$ws.Range("L134").Value = 7790.25
$ws.Range("M134").Value = -4715.667
$ws.Range("N134").Value = -12860.25
$ws.Range("H136").Value = 2663.8333
$ws.Range("I136").Value = 2897.889
$ws.Range("J136").Value = 1961.6666
$ws.Range("K136").Value = 8693.667000000001
$ws.Range("L136").Value = 5884.9998
$ws.Range("M136").Value = -6143.667000000001
$ws.Range("N136").Value = -10984.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1535
$ws.Range("J26").Value = 3000
$ws.Range("L26").Value = 9000
$ws.Range("N26").Value = -9576
$ws.Range("H33").Value = 64
$ws.Range("I33").Value = 40
$ws.Range("J33").Value = 76
$ws.Range("K33").Value = 240
$ws.Range("L33").Value = 456
$ws.Range("M33").Value = 43
$ws.Range("N33").Value = -1022

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4048.9
$ws.Range("I122").Value = 4580.3335
$ws.Range("K122").Value = 13741.0005
$ws.Range("M122").Value = -11291.0005
$ws.Range("H132").Value = 2921.9473
$ws.Range("I132").Value = 2407.4375
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 7222.3125
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -4692.3125
$ws.Range("N132").Value = -22058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3679.5908
$ws.Range("I7").Value = 3643.875
$ws.Range("K7").Value = 3643.875
$ws.Range("M7").Value = -3531.875
$ws.Range("H40").Value = 7208.3335
$ws.Range("I40").Value = 7208.3335
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7208.3335
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7072.3335
$ws.Range("N40").ClearContents()
$ws.Range("H48").Value = 150000
$ws.Range("I48").Value = 150000
$ws.Range("K48").Value = 150000
$ws.Range("M48").Value = -149339
$ws.Range("H122").Value = 6829.6665
$ws.Range("I122").Value = 6489
$ws.Range("K122").Value = 19467
$ws.Range("M122").Value = -17017
$ws.Range("H126").Value = 3679.5908
$ws.Range("I126").Value = 3643.875
$ws.Range("K126").Value = 10931.625
$ws.Range("M126").Value = -8461.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2091.611
$ws.Range("I122").Value = 1790.5625
$ws.Range("K122").Value = 5371.6875
$ws.Range("M122").Value = -2921.6875
$ws.Range("H132").Value = 3484
$ws.Range("I132").Value = 2003.0869
$ws.Range("J132").Value = 11999.25
$ws.Range("K132").Value = 6009.2607
$ws.Range("L132").Value = 35997.75
$ws.Range("M132").Value = -3479.2607
$ws.Range("N132").Value = -41057.75
$ws.Range("H136").Value = 1658.5264
$ws.Range("I136").Value = 1741.5883
$ws.Range("J136").Value = 952.5
$ws.Range("K136").Value = 5224.7649
$ws.Range("L136").Value = 2857.5
$ws.Range("M136").Value = -2674.7649
$ws.Range("N136").Value = -7957.5
